# Generate Report for Handback
# Update the "generated at" / handoff-handback timestamps that are
# refreshed each time the handback status report is (re)generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-29 13:05:35"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
# for the first file row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-29 13:05:30"
$wsZhCn.Range("K2").Value = "2016-08-29 13:05:47"

# de-de sheet: Correspond Handoff Datetime for the first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-29 13:05:55"
